# AHDT1_AH_FAIL.docx template update (Release 1.2.1 / MHD2-159 etc.):
#   - sequencing instrument description updated from
#       "...sequenced on an Illumina NovaSeq 6000 with 150 bp paired end reads."
#     to
#       "...sequenced on an Illumina NovaSeq X Plus (Australian Genome Research
#       Facility) with 150 bp paired end reads."
#   - the cached SAVEDATE field result (bottom of the document) updated to the
#     new save date, "4-Mar-2025".

$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1
$wdFindContinue = 1
$wdReplaceAll   = 2

# 1) Sequencing platform: "NovaSeq 6000" -> "NovaSeq X Plus (Australian Genome
#    Research Facility)". The surrounding "...Illumina " / " with 150 bp..."
#    wording is unchanged.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "NovaSeq 6000",                                          # FindText
    $true,                                                    # MatchCase
    $false,                                                   # MatchWholeWord
    $false,                                                   # MatchWildcards
    $false,                                                   # MatchSoundsLike
    $false,                                                   # MatchAllWordForms
    $true,                                                    # Forward
    $wdFindContinue,                                          # Wrap
    $false,                                                   # Format
    "NovaSeq X Plus (Australian Genome Research Facility)",   # ReplaceWith
    $wdReplaceAll)                                             # Replace

if (-not $found1) {
    throw "Could not find 'NovaSeq 6000' to replace."
}

# 2) Cached SAVEDATE field result text (the field code itself,
#    ` SAVEDATE \@ "d-MMM-yyyy" \* MERGEFORMAT `, is left untouched).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "17-Sep-2024",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    $wdFindContinue,
    $false,
    "4-Mar-2025",
    $wdReplaceAll)

if (-not $found2) {
    throw "Could not find '17-Sep-2024' to replace."
}
